$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.301.12'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.469.43'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '593.38'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '179.03'
$ws.Range('E6').Value = '  +4.42%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.467.88'
$ws.Range('E8').Value = '  -0.49%  '
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('E10').Value = '  +5.50%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.07'
$ws.Range('E11').Value = '  -2.47%  '
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.076.57'
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '32.17'
$ws.Range('E14').Value = '  +11.95%  '
$ws.Range('E15').Value = '  +1.46%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '67.301.37'
$ws.Range('E16').Value = '  +0.76%  '
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.472.15'
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('E19').Value = '  -0.52%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.27'
$ws.Range('E20').Value = '  +1.82%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '389.20'
$ws.Range('E21').Value = '  -0.71%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.85'
$ws.Range('E22').Value = '  -1.11%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '72.75'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.998'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.534'
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '5.71'
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.32'
$ws.Range('E28').Value = '  +1.55%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.176'
$ws.Range('E29').Value = '  -2.91%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.16'
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '23.39'
$ws.Range('E34').Value = '  -0.97%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '7.37'
$ws.Range('E35').Value = '  +0.97%  '
$ws.Range('E37').Value = '  -1.79%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '163.65'
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.72'
$ws.Range('E41').Value = '  +7.11%  '
$ws.Range('E42').Value = '  -0.97%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.836.76'
$ws.Range('E43').Value = '  +1.31%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.60'
$ws.Range('E44').Value = '  -1.09%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '26.03'
$ws.Range('E45').Value = '  -0.32%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0722'
$ws.Range('E46').Value = '  -2.35%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '26.52'
$ws.Range('E47').Value = '  -2.63%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '41.83'
$ws.Range('E48').Value = '  -1.87%  '
$ws.Range('E49').Value = '  -0.91%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '336.66'
$ws.Range('E50').Value = '  +0.26%  '
$ws.Range('E51').Value = '  -2.11%  '
